# Monthly upload data updated:
# Insert a new header row at the top carrying the upload "Month" stamp,
# pushing the existing userId / present_peak_reading / present_off_peak_reading
# table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything down one row so the new Month/"Dec-24" row lands at row 1.
$ws.Rows("1:1").Insert() | Out-Null

# New top row: the month this upload covers.
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = '"Dec-24"'
$ws.Range("B1").NumberFormat = "mmm-yy"

# Keep the same selection Excel would land on after this edit.
$ws.Range("C3").Select() | Out-Null
